# edit.ps1 - applies the two changes described by the target diff:
#
#  1. The table on slide 6 switches from the deck's custom table style
#     ({77092453-C296-492D-A8B9-B2159CE54E94}) to the built-in style
#     {137ADB77-059F-4D0C-8431-113063BBD888}.
#
#  2. The presentation theme's colour scheme is swapped from the
#     "Integral" palette to the stock "Office" palette (dk1, lt1, dk2,
#     lt2, accent1-6, hlink, folHlink) - i.e. theme1.xml's <a:clrScheme>
#     becomes the "Office" colours instead of "Integral".

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{137ADB77-059F-4D0C-8431-113063BBD888}")
    }
}

# --- 2. Theme colour scheme ------------------------------------------
# VBA ThemeColorScheme.Item(n).RGB uses the usual &HBBGGRR long, so the
# RGB hex values below (the stock "Office" theme) are byte-reversed.
$officeColors = @(
    0,          # 1  dk1      000000
    16777215,   # 2  lt1      FFFFFF
    6968388,    # 3  dk2      44546A
    15132391,   # 4  lt2      E7E6E6
    13998939,   # 5  accent1  5B9BD5
    3243501,    # 6  accent2  ED7D31
    10855845,   # 7  accent3  A5A5A5
    49407,      # 8  accent4  FFC000
    12874308,   # 9  accent5  4472C4
    4697456,    # 10 accent6  70AD47
    12673797,   # 11 hlink    0563C1
    7491477     # 12 folHlink 954F72
)

$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
